# "classFields" sheet: the rows listing fields per class were re-ordered
# (e.g. OrderManageService's fields now list SOURCE/template/repository/LOG
# instead of template/repository/SOURCE/LOG, PaymentComponentTests now lists
# template/customer/LOG/kafka/repository instead of
# kafka/repository/template/LOG/customer, and Customer's fields now list
# amountAvailable/id/amountReserved/name instead of
# name/amountAvailable/id/amountReserved).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# pl.piomin.payment.service.OrderManageService
$ws.Cells.Item(2, 2).Value = "SOURCE"
$ws.Cells.Item(2, 4).Value = "java.lang.String"

$ws.Cells.Item(3, 2).Value = "template"
$ws.Cells.Item(3, 4).Value = "org.springframework.kafka.core.KafkaTemplate"

$ws.Cells.Item(4, 2).Value = "repository"
$ws.Cells.Item(4, 4).Value = "pl.piomin.payment.repository.CustomerRepository"

# row 5 (LOG) is unchanged

# pl.piomin.payment.PaymentComponentTests
# row 6 (factory) is unchanged

$ws.Cells.Item(7, 2).Value = "template"
$ws.Cells.Item(7, 4).Value = "org.springframework.kafka.core.KafkaTemplate"

$ws.Cells.Item(8, 2).Value = "customer"
$ws.Cells.Item(8, 4).Value = "pl.piomin.payment.domain.Customer"

$ws.Cells.Item(9, 2).Value = "LOG"
$ws.Cells.Item(9, 4).Value = "org.slf4j.Logger"

$ws.Cells.Item(10, 2).Value = "kafka"
$ws.Cells.Item(10, 4).Value = "org.springframework.kafka.test.EmbeddedKafkaBroker"

$ws.Cells.Item(11, 2).Value = "repository"
$ws.Cells.Item(11, 4).Value = "pl.piomin.payment.repository.CustomerRepository"

# pl.piomin.payment.PaymentApp
# rows 12-13 (LOG, orderManageService) are unchanged

$ws.Cells.Item(14, 2).Value = "repository"
# row 14 column D (pl.piomin.payment.repository.CustomerRepository) is unchanged

# pl.piomin.payment.domain.Customer
$ws.Cells.Item(15, 2).Value = "amountAvailable"
$ws.Cells.Item(15, 4).Value = "int"

$ws.Cells.Item(16, 2).Value = "id"
$ws.Cells.Item(16, 4).Value = "java.lang.Long"

$ws.Cells.Item(17, 2).Value = "amountReserved"
$ws.Cells.Item(17, 4).Value = "int"

$ws.Cells.Item(18, 2).Value = "name"
$ws.Cells.Item(18, 4).Value = "java.lang.String"
